$d = $word.ActiveDocument

# The "Decreases cat's mood" row currently reads "Decreases cat's mood  by 75"
# because the second run starts with a stray leading space (" by 75"). Find
# that run's range first so we know exactly where the edit point is.
$found = $d.Content.Duplicate
$found.Find.Execute(" by 75", $true, $false, $false, $false, $false, $true, `
                     1, $false, "", 0)

# Word stamps the last edit position with the hidden "_GoBack" bookmark.
# Re-adding it here moves it to the new edit location (it is automatically
# removed from its old spot, since bookmark names are unique in a document).
$editPoint = $d.Range($found.Start, $found.Start)
$d.Bookmarks.Add("_GoBack", $editPoint)

# Now clean up the text itself: " by 75" -> "by 75".
$found2 = $d.Content.Duplicate
$found2.Find.Execute(" by 75", $true, $false, $false, $false, $false, $true, `
                      1, $false, "", 0)
$editedRun = $d.Range($found2.Start, $found2.End)
$editedRun.Text = "by 75"
